$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 987.5625
$ws.Range("I11").Value = 987.5625
$ws.Range("K11").Value = 987.5625
$ws.Range("M11").Value = -847.5625
$ws.Range("H18").Value = 2498.3333
$ws.Range("I18").Value = 2498.3333
$ws.Range("K18").Value = 2498.3333
$ws.Range("M18").Value = -2214.3333
$ws.Range("H40").Value = 6539696.5
$ws.Range("I40").Value = 3820.375
$ws.Range("J40").Value = 12349364
$ws.Range("K40").Value = 3820.375
$ws.Range("L40").Value = 12349364
$ws.Range("M40").Value = -3645.375
$ws.Range("N40").Value = -12349714
$ws.Range("H112").Value = 72226.31
$ws.Range("I112").Value = 251849
$ws.Range("J112").Value = 43486.68
$ws.Range("K112").Value = 755547
$ws.Range("L112").Value = 130460.04
$ws.Range("M112").Value = -754439
$ws.Range("N112").Value = -132676.04
$ws.Range("H132").Value = 1750.8462
$ws.Range("I132").Value = 1596.4546
$ws.Range("K132").Value = 4789.3638
$ws.Range("M132").Value = -2259.3638
$ws.Range("H137").Value = 2101.7693
$ws.Range("I137").Value = 1373.625
$ws.Range("J137").Value = 3266.8
$ws.Range("K137").Value = 4120.875
$ws.Range("L137").Value = 9800.400000000001
$ws.Range("M137").Value = -1570.875
$ws.Range("N137").Value = -14900.4
$ws.Range("H138").Value = 4635.2812
$ws.Range("I138").Value = 1219.45
$ws.Range("J138").Value = 6187.9316
$ws.Range("K138").Value = 3658.35
$ws.Range("L138").Value = 18563.7948
$ws.Range("M138").Value = 1481.65
$ws.Range("N138").Value = -28843.7948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1089.1111
$ws.Range("I2").Value = 732.73334
$ws.Range("K2").Value = 732.73334
$ws.Range("M2").Value = -619.73334
$ws.Range("H116").Value = 1089.1111
$ws.Range("I116").Value = 732.73334
$ws.Range("K116").Value = 732.73334
$ws.Range("M116").Value = 1561.26666
$ws.Range("H122").Value = 2667
$ws.Range("I122").Value = 2750.375
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8251.125
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -5801.125
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1089.1111
$ws.Range("I3").Value = 732.73334
$ws.Range("K3").Value = 732.73334
$ws.Range("M3").Value = -618.73334
$ws.Range("H22").Value = 3665722.8
$ws.Range("I22").Value = 1931.3636
$ws.Range("J22").Value = 23816576
$ws.Range("K22").Value = 1931.3636
$ws.Range("L22").Value = 23816576
$ws.Range("M22").Value = -1758.3636
$ws.Range("N22").Value = -23816922
$ws.Range("H80").Value = 353.4375
$ws.Range("I80").Value = 874.6667
$ws.Range("J80").Value = 233.15384
$ws.Range("K80").Value = 874.6667
$ws.Range("L80").Value = 233.15384
$ws.Range("M80").Value = 123.3333
$ws.Range("N80").Value = -2229.15384
$ws.Range("H83").Value = 353.4375
$ws.Range("I83").Value = 874.6667
$ws.Range("J83").Value = 233.15384
$ws.Range("K83").Value = 4373.3335
$ws.Range("L83").Value = 1165.7692
$ws.Range("M83").Value = 618.6665000000003
$ws.Range("N83").Value = -11149.7692
$ws.Range("H94").Value = 312.5
$ws.Range("I94").Value = 312.5
$ws.Range("K94").Value = 312.5
$ws.Range("M94").Value = 138.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 239
$ws.Range("I7").Value = 88
$ws.Range("J7").Value = 346.85715
$ws.Range("K7").Value = 88
$ws.Range("L7").Value = 346.85715
$ws.Range("M7").Value = 25
$ws.Range("N7").Value = -572.85715
$ws.Range("H58").Value = 26322476
$ws.Range("I58").Value = 31257634
$ws.Range("J58").Value = 1635
$ws.Range("K58").Value = 31257634
$ws.Range("L58").Value = 1635
$ws.Range("M58").Value = -31257431
$ws.Range("N58").Value = -2041
$ws.Range("H132").Value = 71431736
$ws.Range("I132").Value = 83336430
$ws.Range("J132").Value = 3599.5
$ws.Range("K132").Value = 250009290
$ws.Range("L132").Value = 10798.5
$ws.Range("M132").Value = -250006760
$ws.Range("N132").Value = -15858.5
$ws.Range("H134").Value = 6251933.5
$ws.Range("I134").Value = 6251933.5
$ws.Range("K134").Value = 18755800.5
$ws.Range("M134").Value = -18753265.5
$ws.Range("H135").Value = 135000
$ws.Range("J135").Value = 135000
$ws.Range("L135").Value = 135000
$ws.Range("N135").Value = -145140
$ws.Range("H136").Value = 26322476
$ws.Range("I136").Value = 31257634
$ws.Range("J136").Value = 1635
$ws.Range("K136").Value = 93772902
$ws.Range("L136").Value = 4905
$ws.Range("M136").Value = -93770352
$ws.Range("N136").Value = -10005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 15375.044
$ws.Range("I56").Value = 15375.044
$ws.Range("K56").Value = 15375.044
$ws.Range("M56").Value = -14845.044
$ws.Range("H107").Value = 1178.75
$ws.Range("I107").Value = 665.8570999999999
$ws.Range("J107").Value = 1577.6666
$ws.Range("K107").Value = 1997.5713
$ws.Range("L107").Value = 4732.9998
$ws.Range("M107").Value = -77.57129999999984
$ws.Range("N107").Value = -8572.9998
$ws.Range("H132").Value = 2079.7727
$ws.Range("I132").Value = 1442.5
$ws.Range("J132").Value = 2221.389
$ws.Range("K132").Value = 12982.5
$ws.Range("L132").Value = 19992.501
$ws.Range("M132").Value = -10452.5
$ws.Range("N132").Value = -25052.501
$ws.Range("H137").Value = 11113099
$ws.Range("I137").Value = 33335334
$ws.Range("K137").Value = 100006002
$ws.Range("M137").Value = -100000902
$ws.Range("H138").Value = 1525.5714
$ws.Range("I138").Value = 1525.5714
$ws.Range("K138").Value = 4576.7142
$ws.Range("M138").Value = 563.2857999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3206.4666
$ws.Range("I80").Value = 3123.5
$ws.Range("K80").Value = 3123.5
$ws.Range("M80").Value = -2125.5
$ws.Range("H83").Value = 3206.4666
$ws.Range("I83").Value = 3123.5
$ws.Range("K83").Value = 15617.5
$ws.Range("M83").Value = -10625.5
$ws.Range("H102").Value = 1062.0555
$ws.Range("I102").Value = 1006.8823
$ws.Range("K102").Value = 1006.8823
$ws.Range("M102").Value = 615.1177
$ws.Range("H132").Value = 9617551
$ws.Range("I132").Value = 10418763
$ws.Range("K132").Value = 31256289
$ws.Range("M132").Value = -31253759

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17867394
$ws.Range("I132").Value = 19241424
$ws.Range("K132").Value = 57724272
$ws.Range("M132").Value = -57721742

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2691.3125
$ws.Range("I96").Value = 3207.2
$ws.Range("J96").Value = 1831.5
$ws.Range("K96").Value = 3207.2
$ws.Range("L96").Value = 1831.5
$ws.Range("M96").Value = -1834.2
$ws.Range("N96").Value = -4577.5
$ws.Range("H132").Value = 27781032
$ws.Range("I132").Value = 45457988
$ws.Range("K132").Value = 136373964
$ws.Range("M132").Value = -136371434
